# Update CDEs (Copy-Paste Fehler bei Codierung behoben),
# neue Informationen unter "Detaillierte Spezifikationen"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-03-04T00:41:34+01:00"

# --- "Include from UCUM" sheet: remove erroneous copy/paste derived-unit rows ---
# Rows (as they exist before any deletion):
#  6  ug/{TotalVolume}
#  7  ug/{Specimen}
#  9  mg/{Volume}
# 10  mg/{TotalVolume}
# 12  g/{TotalWeight}
# Delete from the bottom up so row numbers of not-yet-deleted rows stay valid.
$wsUcum = $wb.Worksheets.Item("Include from UCUM")
$wsUcum.Rows(12).Delete()
$wsUcum.Rows(10).Delete()
$wsUcum.Rows(9).Delete()
$wsUcum.Rows(7).Delete()
$wsUcum.Rows(6).Delete()
